$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-7 and add new rows 8-10 per updated NATMI analysis
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Slit3"
$ws.Range("C2").Value = "Robo4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.044914666666666
$ws.Range("H2").Value = 3.134744
$ws.Range("I2").Value = 0.006668841574421894
$ws.Range("J2").Value = 0.006668841574421893
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 25.33077566666667
$ws.Range("N2").Value = 75.992327
$ws.Range("O2").Value = 0.988229338287255
$ws.Range("P2").Value = 0.988229338287255
$ws.Range("Q2").Value = 26.46849901214311
$ws.Range("R2").Value = 238.216491109288
$ws.Range("S2").Value = 0.006590344896233484
$ws.Range("T2").Value = 0.006590344896233483

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Slit3"
$ws.Range("C3").Value = "Robo4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.044914666666666
$ws.Range("H3").Value = 3.134744
$ws.Range("I3").Value = 0.006668841574421894
$ws.Range("J3").Value = 0.006668841574421893
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05766533333333334
$ws.Range("N3").Value = 0.172996
$ws.Range("O3").Value = 0.002249697164903793
$ws.Range("P3").Value = 0.002249697164903793
$ws.Range("Q3").Value = 0.06025535255822222
$ws.Range("R3").Value = 0.542298173024
$ws.Range("S3").Value = 0.00001500287398316948
$ws.Range("T3").Value = 0.00001500287398316948

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Slit3"
$ws.Range("C4").Value = "Robo4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.044914666666666
$ws.Range("H4").Value = 3.134744
$ws.Range("I4").Value = 0.006668841574421894
$ws.Range("J4").Value = 0.006668841574421893
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.244046
$ws.Range("N4").Value = 0.732138
$ws.Range("O4").Value = 0.009520964547841182
$ws.Range("P4").Value = 0.009520964547841182
$ws.Range("Q4").Value = 0.2550072447413332
$ws.Range("R4").Value = 2.295065202671999
$ws.Range("S4").Value = 0.00006349380420524023
$ws.Range("T4").Value = 0.00006349380420524021

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Slit3"
$ws.Range("C5").Value = "Robo4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 127.867017
$ws.Range("H5").Value = 383.601051
$ws.Range("I5").Value = 0.8160713081836135
$ws.Range("J5").Value = 0.8160713081836134
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 25.33077566666667
$ws.Range("N5").Value = 75.992327
$ws.Range("O5").Value = 0.988229338287255
$ws.Range("P5").Value = 0.988229338287255
$ws.Range("Q5").Value = 3238.970722792853
$ws.Range("R5").Value = 29150.73650513568
$ws.Range("S5").Value = 0.806465608881507
$ws.Range("T5").Value = 0.8064656088815069

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Slit3"
$ws.Range("C6").Value = "Robo4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 127.867017
$ws.Range("H6").Value = 383.601051
$ws.Range("I6").Value = 0.8160713081836135
$ws.Range("J6").Value = 0.8160713081836134
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.05766533333333334
$ws.Range("N6").Value = 0.172996
$ws.Range("O6").Value = 0.002249697164903793
$ws.Range("P6").Value = 0.002249697164903793
$ws.Range("Q6").Value = 7.373494157644
$ws.Range("R6").Value = 66.361447418796
$ws.Range("S6").Value = 0.001835913308380005
$ws.Range("T6").Value = 0.001835913308380005

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Slit3"
$ws.Range("C7").Value = "Robo4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 127.867017
$ws.Range("H7").Value = 383.601051
$ws.Range("I7").Value = 0.8160713081836135
$ws.Range("J7").Value = 0.8160713081836134
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.244046
$ws.Range("N7").Value = 0.732138
$ws.Range("O7").Value = 0.009520964547841182
$ws.Range("P7").Value = 0.009520964547841182
$ws.Range("Q7").Value = 31.20543403078199
$ws.Range("R7").Value = 280.848906277038
$ws.Range("S7").Value = 0.00776978599372656
$ws.Range("T7").Value = 0.007769785993726559

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Slit3"
$ws.Range("C8").Value = "Robo4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 27.77415166666666
$ws.Range("H8").Value = 83.32245499999999
$ws.Range("I8").Value = 0.1772598502419647
$ws.Range("J8").Value = 0.1772598502419647
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 25.33077566666667
$ws.Range("N8").Value = 75.992327
$ws.Range("O8").Value = 0.988229338287255
$ws.Range("P8").Value = 0.988229338287255
$ws.Range("Q8").Value = 703.5408052003095
$ws.Range("R8").Value = 6331.867246802784
$ws.Range("S8").Value = 0.1751733845095147
$ws.Range("T8").Value = 0.1751733845095147

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Slit3"
$ws.Range("C9").Value = "Robo4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 27.77415166666666
$ws.Range("H9").Value = 83.32245499999999
$ws.Range("I9").Value = 0.1772598502419647
$ws.Range("J9").Value = 0.1772598502419647
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.05766533333333334
$ws.Range("N9").Value = 0.172996
$ws.Range("O9").Value = 0.002249697164903793
$ws.Range("P9").Value = 0.002249697164903793
$ws.Range("Q9").Value = 1.601605713908889
$ws.Range("R9").Value = 14.41445142518
$ws.Range("S9").Value = 0.0003987809825406191
$ws.Range("T9").Value = 0.0003987809825406189

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Slit3"
$ws.Range("C10").Value = "Robo4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 27.77415166666666
$ws.Range("H10").Value = 83.32245499999999
$ws.Range("I10").Value = 0.1772598502419647
$ws.Range("J10").Value = 0.1772598502419647
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.244046
$ws.Range("N10").Value = 0.732138
$ws.Range("O10").Value = 0.009520964547841182
$ws.Range("P10").Value = 0.009520964547841182
$ws.Range("Q10").Value = 6.778170617643332
$ws.Range("R10").Value = 61.00353555878999
$ws.Range("S10").Value = 0.001687684749909383
$ws.Range("T10").Value = 0.001687684749909383
